$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1")

# G2: fill in the passenger count, matching the text already used in B2 ("1"),
# via a values-only paste so the existing cell style (s=11) is left untouched.
$ws.Range("B2").Copy()
$ws.Range("G2").PasteSpecial(-4163)  # xlPasteValues

# D8: new empty cell picking up the same style as E9 (s=10).
$ws.Range("E9").Copy()
$ws.Range("D8").PasteSpecial(-4122)  # xlPasteFormats

# Leave the selection on E9, like the final state of the edit.
$ws.Range("E9").Select()
